$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: TestInventory  (D2:E7 item numbers bumped from 90-94 -> 120-124,
# and the very first leather item renumbered 00-18 -> 00-25)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("TestInventory")

$ws1.Range("D2").Value = "Leather00-25"
$ws1.Range("E2").Value = "LT-025"

$ws1.Range("D3").Value = "(Automation) Item 120"
$ws1.Range("E3").Value = "AT-IT-120"

$ws1.Range("D4").Value = "(Automation) Item 121"
$ws1.Range("E4").Value = "AT-IT-121"

$ws1.Range("D5").Value = "(Automation) Item 122"
$ws1.Range("E5").Value = "AT-IT-122"

$ws1.Range("D6").Value = "(Automation) Item 123"
$ws1.Range("E6").Value = "AT-IT-123"

$ws1.Range("D7").Value = "(Automation) Item 124"
$ws1.Range("E7").Value = "AT-IT-124"

$ws1.Activate()
$ws1.Range("F19").Select()

# ---------------------------------------------------------------------------
# Sheet 2: ItemCarousel  (new "leather_flg" + "leather_ball_imgPath" columns
# inserted after h_quantity/imgPath respectively, plus the same item
# renumbering as sheet 1)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("ItemCarousel")

# Insert the two new columns, shifting everything from L onward to the right.
# First insert lands the new "leather_flg" column at L (old L/M/N -> M/N/O).
$ws2.Range("L1").EntireColumn.Insert()
# Second insert lands "leather_ball_imgPath" at N (old M/N, now at N/O -> O/P).
$ws2.Range("N1").EntireColumn.Insert()

# Header row
$ws2.Range("L1").Value = "leather_flg"
$ws2.Range("N1").Value = "leather_ball_imgPath"

# Row 2
$ws2.Range("C2").Value = "Leather00-24"
$ws2.Range("D2").Value = "LT-024"
$ws2.Range("E2").Value = "T0-LT-00-24"
$ws2.Range("F2").Value = "AT-LT-00-24"
$ws2.Range("G2").Value = "Auto Testing -013"
$ws2.Range("L2").Value = 1
$ws2.Range("N2").Value = "E:\POM_for_Bottle\WebApp\TestData\InventoryData\Image\David Leather Jacket.jpg"

# Row 3
$ws2.Range("C3").Value = "(Automation) Item 120"
$ws2.Range("D3").Value = "AT-IT-120"
$ws2.Range("E3").Value = "T0120"
$ws2.Range("F3").Value = "AT0120"
$ws2.Range("G3").Value = "Auto Testing 120"
$ws2.Range("L3").Value = 0
$ws2.Range("N3").ClearContents()

# Row 4
$ws2.Range("C4").Value = "(Automation) Item 121"
$ws2.Range("D4").Value = "AT-IT-121"
$ws2.Range("E4").Value = "T0121"
$ws2.Range("F4").Value = "AT0121"
$ws2.Range("G4").Value = "Auto Testing 121"
$ws2.Range("L4").Value = 0
$ws2.Range("N4").ClearContents()

# Row 5
$ws2.Range("C5").Value = "(Automation) Item 122"
$ws2.Range("D5").Value = "AT-IT-122"
$ws2.Range("E5").Value = "T0122"
$ws2.Range("F5").Value = "AT0122"
$ws2.Range("G5").Value = "Auto Testing 122"
$ws2.Range("L5").Value = 0
$ws2.Range("N5").ClearContents()

# Row 6
$ws2.Range("C6").Value = "(Automation) Item 123"
$ws2.Range("D6").Value = "AT-IT-123"
$ws2.Range("E6").Value = "T0123"
$ws2.Range("F6").Value = "AT0123"
$ws2.Range("G6").Value = "Auto Testing 123"
$ws2.Range("L6").Value = 0
$ws2.Range("N6").ClearContents()

# Row 7
$ws2.Range("C7").Value = "(Automation) Item 124"
$ws2.Range("D7").Value = "AT-IT-124"
$ws2.Range("E7").Value = "T0124"
$ws2.Range("F7").Value = "AT0124"
$ws2.Range("G7").Value = "Auto Testing 124"
$ws2.Range("L7").Value = 0
$ws2.Range("N7").ClearContents()

$ws2.Activate()
$ws2.Range("H20").Select()

# ---------------------------------------------------------------------------
# Sheet 3: restocking  (C2:D7 item numbers, same renumbering as sheet 1/2)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("restocking")

$ws3.Range("C2").Value = "Leather00-24"
$ws3.Range("D2").Value = "LT-024"

$ws3.Range("C3").Value = "(Automation) Item 120"
$ws3.Range("D3").Value = "AT-IT-120"

$ws3.Range("C4").Value = "(Automation) Item 121"
$ws3.Range("D4").Value = "AT-IT-121"

$ws3.Range("C5").Value = "(Automation) Item 122"
$ws3.Range("D5").Value = "AT-IT-122"

$ws3.Range("C6").Value = "(Automation) Item 123"
$ws3.Range("D6").Value = "AT-IT-123"

$ws3.Range("C7").Value = "(Automation) Item 124"
$ws3.Range("D7").Value = "AT-IT-124"

$ws3.Activate()
$ws3.Range("C3").Select()

# ---------------------------------------------------------------------------
# Sheet 4: searchItem  (selection only; keep it the last activated sheet so
# it remains the active tab, matching the original workbook state)
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("searchItem")
$ws4.Activate()
$ws4.Range("C3").Select()
